# Updates cryptos list values (Price / Volume(1h) columns) per latest data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.313.75"
$ws.Range("E2").Value = "  +0.54%  "
$ws.Range("D3").Value = "1.667.00"
$ws.Range("E4").Value = "  +0.28%  "
$ws.Range("D5").Value = "'220.60"
$ws.Range("E5").Value = "  +1.37%  "
$ws.Range("D6").Value = "'0.5311"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  +0.26%  "
$ws.Range("D8").Value = "'0.2648"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("D9").Value = "'0.06366"
$ws.Range("E9").Value = "  +0.58%  "
$ws.Range("D10").Value = "'20.85"
$ws.Range("E10").Value = "  +2.43%  "
$ws.Range("D11").Value = "'0.07833"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "'4.513"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Value = "1.670.93"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").Value = "1.895.43"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").Value = "'0.5592"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "0.0₅8162"
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("D17").Value = "'65.81"
$ws.Range("E17").Value = "  +0.69%  "
$ws.Range("D18").Value = "26.312.70"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "'4.714"
$ws.Range("E20").Value = "  +2.75%  "
$ws.Range("D21").Value = "'196.86"
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("D22").Value = "'10.25"
$ws.Range("E22").Value = "  +1.71%  "
$ws.Range("D23").Value = "'6.047"
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").Value = "'1.010"
$ws.Range("E24").Value = "  +0.33%  "
$ws.Range("D25").Value = "'145.87"
$ws.Range("E25").Value = "  +0.49%  "
$ws.Range("D26").Value = "'0.1220"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").Value = "'7.237"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").Value = "'1.507"
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("D30").Value = "'0.05892"
$ws.Range("D31").Value = "'1.284"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").Value = "'3.537"
$ws.Range("D33").Value = "'3.333"
$ws.Range("E33").Value = "  +2.17%  "
$ws.Range("D34").Value = "'1.602"
$ws.Range("E34").Value = "  +1.08%  "
$ws.Range("D35").Value = "'2.831"
$ws.Range("E35").Value = "  +0.83%  "
$ws.Range("D36").Value = "'0.9596"
$ws.Range("E36").Value = "  +1.22%  "
$ws.Range("D37").Value = "'2.435"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").Value = "'0.5815"
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("D40").Value = "'5.945"
$ws.Range("D41").Value = "1.077.61"
$ws.Range("E41").Value = "  +3.63%  "
$ws.Range("D42").Value = "'0.8590"
$ws.Range("E42").Value = "  +1.15%  "
$ws.Range("D43").Value = "'1.009"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").Value = "'102.70"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").Value = "1.806.03"
$ws.Range("D46").Value = "'58.40"
$ws.Range("E46").Value = "  +2.82%  "
$ws.Range("D47").Value = "0.0₈106"
$ws.Range("E47").Value = "  +1.21%  "
$ws.Range("D48").Value = "'1.015"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").Value = "'0.4409"
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("D50").Value = "'8.000"
$ws.Range("E50").Value = "  +1.76%  "
$ws.Range("D51").Value = "'0.05151"
